$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "reftype" column after sequenceID (new column B) ---
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "reftype"
$ws.Range("B2").Value = "Consensus"
$ws.Range("B3").Value = "Consensus"
$ws.Range("B4").Value = "Consensus"
$ws.Range("B5").Value = "Consensus"
$ws.Range("B6").Value = "Consensus"
$ws.Range("B7").Value = "Consensus"
$ws.Range("B8").Value = "Consensus"
$ws.Range("B9").Value = "Consensus"
$ws.Range("B10").Value = "Consensus"
$ws.Range("B11").Value = "Consensus"
$ws.Range("B12").Value = "Consensus"
$ws.Range("B13").Value = "Consensus"
$ws.Range("B14").Value = "Consensus"

# --- Insert 5 new host/locus columns after clade (new columns H:L) ---
$ws.Range("H1:L1").EntireColumn.Insert()

$ws.Range("H1").Value = "host_group_taxlevel"
$ws.Range("I1").Value = "host_group_name"
$ws.Range("J1").Value = "host_group_common_name"
$ws.Range("K1").Value = "num_copies"
$ws.Range("L1").Value = "locus_id"

$ws.Range("H2:L14").Value = 1

# match the plain (unfilled) number style used elsewhere in the sheet
$ws.Range("F7").Copy()
$ws.Range("H2:L14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths (best effort) ---
$ws.Columns.Item(1).ColumnWidth = 41.6640625
$ws.Columns.Item(2).ColumnWidth = 15.1640625
$ws.Columns.Item(3).ColumnWidth = 35
$ws.Columns.Item(4).ColumnWidth = 63
$ws.Columns.Item(5).ColumnWidth = 15.5
$ws.Columns.Item(6).ColumnWidth = 13.83203125
$ws.Columns.Item(7).ColumnWidth = 9.6640625
$ws.Columns.Item(8).ColumnWidth = 24.1640625
$ws.Columns.Item(9).ColumnWidth = 22.1640625
$ws.Columns.Item(10).ColumnWidth = 15.83203125

# --- Selection ---
$ws.Range("A8").Select()
